$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Data edit: Food budget for Jan (C9) drops from 300 to 200 ---
# All dependent formulas (row totals, percentages, column sum/average)
# recalculate automatically.
$ws.Range("C9").Value = 200

# --- View state: scroll the window up/left and move the selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("A5").Select()

# --- Conditional formatting: flag monthly bill amounts over 250 in red ---
$cfRange = $ws.Range("C5:F10")
$cfRange.FormatConditions.Delete()
$condition = $cfRange.FormatConditions.Add(1, 5, "250")  # xlCellValue, xlGreater
$condition.Interior.Color = 255  # RGB(255,0,0) red fill
